# Loading time update after build success from jenkins
# Adds a new "build" column (W) to the LoadTime sheet, mirroring the
# existing " Oct 19" header already used by the previous columns and
# filling in the newly measured loading times for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoadTime")

$ws.Range("W1").Value = $ws.Range("V1").Value

$ws.Range("W2").Value = 0
$ws.Range("W3").Value = 5
$ws.Range("W4").Value = 5
$ws.Range("W5").Value = 0
$ws.Range("W6").Value = 19
$ws.Range("W7").Value = 0
